# Apply row-content rotation as described by the diff.
# Rows 6<->7, 8->9->10->8(cycle), 11->12->13->11(cycle... see mapping), 14->15->16->14(cycle), 23<->24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: pull in the data previously on row 7
$ws.Range("A6").Value = 131046755
$ws.Range("B6").Value = 57881
$ws.Range("E6").Value = 100049
$ws.Range("F6").Value = "Spillkråka"
$ws.Range("G6").Value = "Dryocopus martius"
$ws.Range("H6").Value = "(Linnaeus, 1758)"
$ws.Range("M6").Value = "färska spår"
$ws.Range("Q6").Value = 402424
$ws.Range("R6").Value = 6818357
$ws.Range("Z6").Value = "16:56"
$ws.Range("AB6").Value = "16:56"

# Row 7: pull in the data previously on row 6
$ws.Range("A7").Value = 131046733
$ws.Range("B7").Value = 91808
$ws.Range("E7").Value = 1202
$ws.Range("F7").Value = "Ullticka"
$ws.Range("G7").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H7").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M7").ClearContents()
$ws.Range("Q7").Value = 402493
$ws.Range("R7").Value = 6818443
$ws.Range("Z7").Value = "16:43"
$ws.Range("AB7").Value = "16:43"

# Row 8: pull in the data previously on row 9
$ws.Range("A8").Value = 131046844
$ws.Range("B8").Value = 79243
$ws.Range("E8").Value = 6425
$ws.Range("F8").Value = "Garnlav"
$ws.Range("G8").Value = "Alectoria sarmentosa"
$ws.Range("H8").Value = "(Ach.) Ach."
$ws.Range("Q8").Value = 402484
$ws.Range("R8").Value = 6818538
$ws.Range("Z8").Value = "16:23"
$ws.Range("AB8").Value = "16:23"

# Row 9: pull in the data previously on row 10
$ws.Range("A9").Value = 131046843
$ws.Range("Q9").Value = 402432
$ws.Range("R9").Value = 6818480
$ws.Range("Z9").Value = "16:20"
$ws.Range("AB9").Value = "16:20"

# Row 10: pull in the data previously on row 8
$ws.Range("A10").Value = 131046711
$ws.Range("B10").Value = 83223
$ws.Range("E10").Value = 6440
$ws.Range("F10").Value = "Vitgrynig nållav"
$ws.Range("G10").Value = "Chaenotheca subroscida"
$ws.Range("H10").Value = "(Eitner) Zahlbr."
$ws.Range("Q10").Value = 402363
$ws.Range("R10").Value = 6818428
$ws.Range("Z10").Value = "16:09"
$ws.Range("AB10").Value = "16:09"

# Row 11: pull in the data previously on row 13
$ws.Range("A11").Value = 131046788
$ws.Range("M11").Value = "färska spår"
$ws.Range("Q11").Value = 402473
$ws.Range("R11").Value = 6818425
$ws.Range("Z11").Value = "16:47"
$ws.Range("AB11").Value = "16:47"
$ws.Range("AC11").Value = "Färska ringhack (gran)"

# Row 12: pull in the data previously on row 11
$ws.Range("A12").Value = 131046735
$ws.Range("B12").Value = 57884
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 100109
$ws.Range("F12").Value = "Tretåig hackspett"
$ws.Range("G12").Value = "Picoides tridactylus"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "nyligen använt bo"
$ws.Range("Q12").Value = 402448
$ws.Range("R12").Value = 6818295
$ws.Range("Z12").Value = "16:54"
$ws.Range("AB12").Value = "16:54"

# Row 13: pull in the data previously on row 12
$ws.Range("A13").Value = 131046763
$ws.Range("B13").Value = 92267
$ws.Range("D13").Value = "VU"
$ws.Range("E13").Value = 1209
$ws.Range("F13").Value = "Rynkskinn"
$ws.Range("G13").Value = "Hermanssonia centrifuga"
$ws.Range("H13").Value = "(P. Karst.) Zmitr."
$ws.Range("M13").ClearContents()
$ws.Range("Q13").Value = 402378
$ws.Range("R13").Value = 6818392
$ws.Range("Z13").Value = "17:02"
$ws.Range("AB13").Value = "17:02"
$ws.Range("AC13").ClearContents()

# Row 14: pull in the data previously on row 16
$ws.Range("A14").Value = 131046806
$ws.Range("B14").Value = 83206
$ws.Range("D14").Value = "LC"
$ws.Range("E14").Value = 6439
$ws.Range("F14").Value = "Gulnål"
$ws.Range("G14").Value = "Chaenotheca brachypoda"
$ws.Range("H14").Value = "(Ach.) Tibell"
$ws.Range("Q14").Value = 402340
$ws.Range("R14").Value = 6818363
$ws.Range("Z14").Value = "17:05"
$ws.Range("AB14").Value = "17:05"

# Row 15: pull in the data previously on row 14
$ws.Range("A15").Value = 131046811
$ws.Range("Q15").Value = 402450
$ws.Range("R15").Value = 6818298
$ws.Range("Z15").Value = "16:54"
$ws.Range("AB15").Value = "16:54"

# Row 16: pull in the data previously on row 15
$ws.Range("A16").Value = 131046808
$ws.Range("B16").Value = 91828
$ws.Range("D16").Value = "NT"
$ws.Range("E16").Value = 5432
$ws.Range("F16").Value = "Granticka"
$ws.Range("G16").Value = "Porodaedalea chrysoloma s.lat."
$ws.Range("H16").Value = ""
$ws.Range("Q16").Value = 402323
$ws.Range("R16").Value = 6818416
$ws.Range("Z16").Value = "16:06"
$ws.Range("AB16").Value = "16:06"

# Row 23: pull in the data previously on row 24
$ws.Range("A23").Value = 131046845
$ws.Range("B23").Value = 79243
$ws.Range("E23").Value = 6425
$ws.Range("F23").Value = "Garnlav"
$ws.Range("G23").Value = "Alectoria sarmentosa"
$ws.Range("H23").Value = "(Ach.) Ach."
$ws.Range("M23").ClearContents()
$ws.Range("Q23").Value = 402575
$ws.Range("R23").Value = 6818545
$ws.Range("Z23").Value = "16:34"
$ws.Range("AB23").Value = "16:34"
$ws.Range("AC23").ClearContents()
$ws.Range("AE23").Value = $false

# Row 24: pull in the data previously on row 23
$ws.Range("A24").Value = 131047016
$ws.Range("B24").Value = 57884
$ws.Range("E24").Value = 100109
$ws.Range("F24").Value = "Tretåig hackspett"
$ws.Range("G24").Value = "Picoides tridactylus"
$ws.Range("H24").Value = "(Linnaeus, 1758)"
$ws.Range("M24").Value = "färska spår"
$ws.Range("Q24").Value = 402474
$ws.Range("R24").Value = 6818507
$ws.Range("Z24").Value = "16:22"
$ws.Range("AB24").Value = "16:22"
$ws.Range("AC24").Value = "Troliga spår efter tretåig hackspett (barkfälkning)"
$ws.Range("AE24").Value = $true
